$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 38462196
$ws.Range("I33").Value = 297.10526
$ws.Range("J33").Value = 142858780
$ws.Range("K33").Value = 297.10526
$ws.Range("L33").Value = 142858780
$ws.Range("M33").Value = -68.10525999999999
$ws.Range("N33").Value = -142859238
$ws.Range("H76").Value = 3133.3333
$ws.Range("I76").Value = 3118.182
$ws.Range("K76").Value = 3118.182
$ws.Range("M76").Value = -2803.182
$ws.Range("H79").Value = 3133.3333
$ws.Range("I79").Value = 3118.182
$ws.Range("K79").Value = 3118.182
$ws.Range("M79").Value = -2026.182
$ws.Range("H88").Value = 9063.875
$ws.Range("I88").Value = 8501.5
$ws.Range("J88").Value = 9251.333000000001
$ws.Range("K88").Value = 8501.5
$ws.Range("L88").Value = 9251.333000000001
$ws.Range("M88").Value = -8095.5
$ws.Range("N88").Value = -10063.333
$ws.Range("H91").Value = 9063.875
$ws.Range("I91").Value = 8501.5
$ws.Range("J91").Value = 9251.333000000001
$ws.Range("K91").Value = 8501.5
$ws.Range("L91").Value = 9251.333000000001
$ws.Range("M91").Value = -7097.5
$ws.Range("N91").Value = -12059.333
$ws.Range("H132").Value = 1743.9344
$ws.Range("I132").Value = 882.2941
$ws.Range("J132").Value = 6138.3
$ws.Range("K132").Value = 2646.8823
$ws.Range("L132").Value = 18414.9
$ws.Range("M132").Value = -116.8822999999998
$ws.Range("N132").Value = -23474.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1742.3043
$ws.Range("I132").Value = 1409.0754
$ws.Range("J132").Value = 2846.125
$ws.Range("K132").Value = 4227.2262
$ws.Range("L132").Value = 8538.375
$ws.Range("M132").Value = -1697.2262
$ws.Range("N132").Value = -13598.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1711.5385
$ws.Range("I86").Value = 1588.8889
$ws.Range("J86").Value = 1987.5
$ws.Range("K86").Value = 1588.8889
$ws.Range("L86").Value = 1987.5
$ws.Range("M86").Value = -465.8888999999999
$ws.Range("N86").Value = -4233.5
$ws.Range("H89").Value = 1711.5385
$ws.Range("I89").Value = 1588.8889
$ws.Range("J89").Value = 1987.5
$ws.Range("K89").Value = 7944.4445
$ws.Range("L89").Value = 9937.5
$ws.Range("M89").Value = -2328.4445
$ws.Range("N89").Value = -21169.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19609296
$ws.Range("I31").Value = 33334104
$ws.Range("J31").Value = 2424.8096
$ws.Range("K31").Value = 33334104
$ws.Range("L31").Value = 2424.8096
$ws.Range("M31").Value = -33333809
$ws.Range("N31").Value = -3014.8096
$ws.Range("H34").Value = 19609296
$ws.Range("I34").Value = 33334104
$ws.Range("J34").Value = 2424.8096
$ws.Range("K34").Value = 33334104
$ws.Range("L34").Value = 2424.8096
$ws.Range("M34").Value = -33333902
$ws.Range("N34").Value = -2828.8096
$ws.Range("H99").Value = 3312
$ws.Range("I99").Value = 2925.3
$ws.Range("J99").Value = 4037.0625
$ws.Range("K99").Value = 2925.3
$ws.Range("L99").Value = 4037.0625
$ws.Range("M99").Value = -1427.3
$ws.Range("N99").Value = -7033.0625
$ws.Range("H126").Value = 3312
$ws.Range("I126").Value = 2925.3
$ws.Range("J126").Value = 4037.0625
$ws.Range("K126").Value = 8775.900000000001
$ws.Range("L126").Value = 12111.1875
$ws.Range("M126").Value = -6305.900000000001
$ws.Range("N126").Value = -17051.1875
$ws.Range("H132").Value = 2643.9666
$ws.Range("I132").Value = 2119.1365
$ws.Range("J132").Value = 4087.25
$ws.Range("K132").Value = 6357.4095
$ws.Range("L132").Value = 12261.75
$ws.Range("M132").Value = -3827.4095
$ws.Range("N132").Value = -17321.75
$ws.Range("H134").Value = 2161.4814
$ws.Range("I134").Value = 2110
$ws.Range("J134").Value = 3500
$ws.Range("K134").Value = 6330
$ws.Range("L134").Value = 10500
$ws.Range("M134").Value = -3795
$ws.Range("N134").Value = -15570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 3651.25
$ws.Range("I17").Value = 986
$ws.Range("J17").Value = 5250.4
$ws.Range("K17").Value = 2958
$ws.Range("L17").Value = 15751.2
$ws.Range("M17").Value = -2789
$ws.Range("N17").Value = -16089.2
$ws.Range("H68").Value = 768.0909
$ws.Range("I68").Value = 616.6667
$ws.Range("J68").Value = 949.8
$ws.Range("K68").Value = 1850.0001
$ws.Range("L68").Value = 2849.4
$ws.Range("M68").Value = -1039.0001
$ws.Range("N68").Value = -4471.4
$ws.Range("H71").Value = 768.0909
$ws.Range("I71").Value = 616.6667
$ws.Range("J71").Value = 949.8
$ws.Range("K71").Value = 5550.0003
$ws.Range("L71").Value = 8548.199999999999
$ws.Range("M71").Value = -1494.0003
$ws.Range("N71").Value = -16660.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5249.8
$ws.Range("I70").Value = 5090.476
$ws.Range("J70").Value = 5488.7856
$ws.Range("K70").Value = 5090.476
$ws.Range("L70").Value = 5488.7856
$ws.Range("M70").Value = -4820.476
$ws.Range("N70").Value = -6028.7856
$ws.Range("H73").Value = 5249.8
$ws.Range("I73").Value = 5090.476
$ws.Range("J73").Value = 5488.7856
$ws.Range("K73").Value = 5090.476
$ws.Range("L73").Value = 5488.7856
$ws.Range("M73").Value = -4154.476
$ws.Range("N73").Value = -7360.7856
$ws.Range("H102").Value = 1546.1111
$ws.Range("I102").Value = 1315.5385
$ws.Range("K102").Value = 1315.5385
$ws.Range("M102").Value = 306.4614999999999
$ws.Range("H132").Value = 2659.5
$ws.Range("I132").Value = 2783.6428
$ws.Range("J132").Value = 2311.9
$ws.Range("K132").Value = 8350.928400000001
$ws.Range("L132").Value = 6935.700000000001
$ws.Range("M132").Value = -5820.928400000001
$ws.Range("N132").Value = -11995.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 23178.2
$ws.Range("I68").Value = 35730.332
$ws.Range("J68").Value = 4350
$ws.Range("K68").Value = 35730.332
$ws.Range("L68").Value = 4350
$ws.Range("M68").Value = -34981.332
$ws.Range("N68").Value = -5848
$ws.Range("H71").Value = 23178.2
$ws.Range("I71").Value = 35730.332
$ws.Range("J71").Value = 4350
$ws.Range("K71").Value = 178651.66
$ws.Range("L71").Value = 21750
$ws.Range("M71").Value = -174907.66
$ws.Range("N71").Value = -29238

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3574009.8
$ws.Range("I132").Value = 4350262
$ws.Range("J132").Value = 3249.8
$ws.Range("K132").Value = 13050786
$ws.Range("L132").Value = 9749.400000000001
$ws.Range("M132").Value = -13048256
$ws.Range("N132").Value = -14809.4
$ws.Range("H136").Value = 20001860
$ws.Range("I136").Value = 25001808
$ws.Range("J136").Value = 2069.8
$ws.Range("K136").Value = 75005424
$ws.Range("L136").Value = 6209.400000000001
$ws.Range("M136").Value = -75002874
$ws.Range("N136").Value = -11309.4
